# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> currently the plain "Office Theme" palette,
#                            used only by the Notes Master.
#   ppt/theme/theme2.xml -> currently the "Integral" palette, used by the
#                            (single) Slide Master / the whole deck design.
#
# The authored change swaps the two themes' contents: the deck's live
# design (theme2.xml) becomes the plain "Office Theme" colors, and the
# notes-only theme (theme1.xml) becomes the "Integral" colors. The
# fontScheme/fmtScheme bodies are identical between the two themes, so the
# only real content difference is the 12 color-scheme entries (and the
# theme/colorScheme display names, which PowerPoint does not allow
# re-pointing through the exposed object model).
#
# Apply the reachable half of that swap: push the "Office Theme" RGB
# values onto the presentation's live theme color scheme (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) through the Slide Master's Theme object.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Index order for ThemeColorScheme.Item(n): 1=dk1 2=lt1 3=dk2 4=lt2
# 5=accent1 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink
# 12=folHlink. RGB() isn't available in this host, so values are supplied
# as the packed 0xBBGGRR integer PowerPoint's RGB property expects.

$colorScheme.Item(1).RGB = 0          # dk1      = 000000
$colorScheme.Item(2).RGB = 16777215   # lt1      = FFFFFF
$colorScheme.Item(3).RGB = 6968388    # dk2      = 44546A
$colorScheme.Item(4).RGB = 15132391   # lt2      = E7E6E6
$colorScheme.Item(5).RGB = 13998939   # accent1  = 5B9BD5
$colorScheme.Item(6).RGB = 3243501    # accent2  = ED7D31
$colorScheme.Item(7).RGB = 10855845   # accent3  = A5A5A5
$colorScheme.Item(8).RGB = 49407      # accent4  = FFC000
$colorScheme.Item(9).RGB = 12874308   # accent5  = 4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  = 70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    = 0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink = 954F72
